# Update countries & provincias Spain
# Applies the 29-Jul-2020 20:48 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 20:48"

# --- Straightforward per-country numeric refreshes (no re-sort needed) -
# Columns: A Pais | B Casos totales | C Nuevos casos | D Casos activos |
#          E Recuperados | F Casos criticos | G Muertes hoy | H Muertes

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4529658
$ws.Range("C4").Value = 31315
$ws.Range("D4").Value = 2212272
$ws.Range("E4").Value = 2164391
$ws.Range("G4").Value = 704
$ws.Range("H4").Value = 152995

# Row 6 - India
$ws.Range("B6").Value = 1584384
$ws.Range("C6").Value = 52249
$ws.Range("D6").Value = 1021611
$ws.Range("E6").Value = 527770
$ws.Range("G6").Value = 779
$ws.Range("H6").Value = 35003

# Row 21 - Alemania
$ws.Range("B21").Value = 208666
$ws.Range("C21").Value = 715
$ws.Range("E21").Value = 7455

# Row 22 - Francia
$ws.Range("B22").Value = 185196
$ws.Range("C22").Value = 1392
$ws.Range("E22").Value = 73647
$ws.Range("G22").Value = 15
$ws.Range("H22").Value = 30238

# Row 25 - Canada
$ws.Range("E25").Value = 6024
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 8914

# Row 48 - Guatemala
$ws.Range("B48").Value = 47605
$ws.Range("C48").Value = 1154
$ws.Range("D48").Value = 34488
$ws.Range("E48").Value = 11282
$ws.Range("G48").Value = 53
$ws.Range("H48").Value = 1835

# Row 108 - Maldivas
$ws.Range("B108").Value = 3567
$ws.Range("C108").Value = 61
$ws.Range("D108").Value = 2554
$ws.Range("E108").Value = 998

# Row 122 - Cabo Verde
$ws.Range("B122").Value = 2373
$ws.Range("C122").Value = 19
$ws.Range("D122").Value = 1694
$ws.Range("E122").Value = 656
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 23

# Row 123 - Sudan del Sur
$ws.Range("B123").Value = 2322
$ws.Range("C123").Value = 17
$ws.Range("E123").Value = 1101

# --- Siria overtakes Crucero: rows 156-158 re-sort -----------------------
# Before: 156 Crucero, 157 San Marino, 158 Siria
# After:  156 Siria (updated figures), 157 Crucero, 158 San Marino
# Crucero's and San Marino's own figures are unchanged, they just shift
# down one row to make room for Siria's improved (higher) total.

$ws.Range("A156").Value = "Siria"
$ws.Range("B156").Value = 717
$ws.Range("C156").Value = 23
$ws.Range("D156").Value = 229
$ws.Range("E156").Value = 448
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 40

$ws.Range("A157").Value = "Crucero"
$ws.Range("B157").Value = 712
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 651
$ws.Range("E157").Value = 48
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 13

$ws.Range("A158").Value = "San Marino"
$ws.Range("B158").Value = 699
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 657
$ws.Range("E158").Value = 0
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 42
